# Generate Report for Handoff
# The 388cb62e-4b38-4fce-8e47-1c7b6e786d71 file has left the batch; the
# 131a7ed5-056e-4060-bc28-41101af30063 file has been (re)handed off instead
# of handed back, with fresh timestamps. So: drop each sheet's row 3
# (the 388cb62e entry) and update row 2's status/timestamp cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Drop every hyperlink on the sheet (this runtime's Hyperlinks.Delete()
# is sheet-wide regardless of the Range it's called through), then
# delete row 3 and re-add the hyperlink(s) that must survive.
$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-22 12:55:47"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/38baa5749e747bcb78e69737433e7f3dac72c681/e2e/131a7ed5-056e-4060-bc28-41101af30063.md", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-22 12:55:43"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/38baa5749e747bcb78e69737433e7f3dac72c681/e2e/131a7ed5-056e-4060-bc28-41101af30063.md", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9c84f99e9046e2999291436e7c13968fb1a4f0e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/36f8025966ba4b48b3ac36ce2ccba273db63720a/e2e/131a7ed5-056e-4060-bc28-41101af30063.md", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a3c554dc4d73e86f8629d5537cd73269ecb3124e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-22 12:55:47"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/38baa5749e747bcb78e69737433e7f3dac72c681/e2e/131a7ed5-056e-4060-bc28-41101af30063.md", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1afc6814547d81ca35b2a608210821a33fc8d556/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/214df847fee3b63dd7abcc8c102dd647abf481e8/e2e/131a7ed5-056e-4060-bc28-41101af30063.md", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/41c14018991f20d24f2e9c30fc15b4e1ff49d371/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf", [Type]::Missing, [Type]::Missing, "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf") | Out-Null

$wb.Save()
